$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the full content of data rows 2 and 3 (everything except the
# header row stays put). Only the columns that actually carry data in either
# row are touched, so we don't create spurious empty cells in unused columns.
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","N","P","Q","R","S","T","U","V","W","Y","Z","AA","AB","AD","AE","AF","AG","AT","AW","AX","AY")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")
    $v2 = $cell2.Value2
    $v3 = $cell3.Value2
    $cell2.Value = $v3
    $cell3.Value = $v2
}

# Row 3 originally had an explicit (empty) cell in column L while row 2 did
# not; after the swap row 2 should carry that placeholder cell and row 3
# should no longer have one. Value = "" clears/removes cell content, so we
# mark the now-empty L2 cell as present the same way the source workbook
# represented it - via an explicit (blank) number format on the cell.
$ws.Range("L2").NumberFormat = $ws.Range("L2").NumberFormat
$ws.Range("L3").ClearContents()
